# Generate Report for Handoff
# Update the "Latest Handoff Date"/"Latest Handoff Datetime" entries for the
# 52679a15-f857-4c26-9a01-c476e58b1a39 file (row 6 in every table) to reflect
# a freshly generated handoff report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest Handoff Date" column (D) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D6").Value = "2016-34-13 00:34:37"

# --- zh-cn sheet: "Latest Handoff Datetime" column (E) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E6").Value = "2016-03-13 00:34:34"

# --- de-de sheet: "Latest Handoff Datetime" column (E) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E6").Value = "2016-03-13 00:34:37"
